$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift header row: insert two new "week" columns (Jun_15, Jun_17) before the
# existing Jun_13 / Jun_10 columns, moving data from C->D and B->C conceptually,
# then filling in the new leftmost columns.
$ws.Range("E1").Value = "Jun_10"
$ws.Range("D1").Value = "Jun_13"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Fill the two new data columns (D and E) for every data row with the same
# "UN" marker used by column C, except row 8 which carries the long
# analyst-action note in the new last column (E), matching column C's
# special value for that row.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 4).Value = "UN"
    $ws.Cells.Item($r, 5).Value = "UN"
}

# Row 8's special analyst-action note used to live in column C; it now moves
# to the new last column (E), and column C reverts to the normal "UN" marker.
$ws.Range("C8").Value = "UN"
$ws.Range("E8").Value = "12/5/2017,Initiated Coverage,Sector Weight ➝ Sector Weight,"

# Column widths: keep column C at width 8 (now marked collapsed) and give the
# two newly inserted columns (D, E) the same width.
$ws.Columns("C").ColumnWidth = 7.16
$ws.Columns("D").ColumnWidth = 7.16
$ws.Columns("E").ColumnWidth = 7.16
